# Añadido tipo de pago
# Insert two new columns (H:I) for "Tipo de pago" (payment type) and
# "Valor Pagado" (amount paid), pushing the old "Valor" column (and
# everything after it) two columns to the right.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift columns H..P right by two, opening up H:I for the new fields.
$ws.Range("H1:I1").EntireColumn.Insert()

# New header cells.
$ws.Range("H1").Value = "Tipo de pago"
$ws.Range("I1").Value = "Valor Pagado"

# New data cells.
$ws.Range("H2").Value = 1
$ws.Range("I2").Value = 1500000

$ws.Range("H3").Value = 2
$ws.Range("I3").Value = 120000

# The old "Valor" column (now shifted to column J) gets new amounts too.
$ws.Range("J2").Value = 1500000
$ws.Range("J3").Value = 1500000
